# Update the "苏州-漫展信息" workbook to the scraper output generated at
# a56beed: refresh the "想去人数" counters on the existing exhibition rows
# and insert a new entry - 2024.03.08 "苏州·国风宠物-cosplay展" - as the new
# row 12, pushing the later rows down by one. Applies to both the "展览"
# and "全部类型" sheets, which hold identical copies of the table.

function Set-TextValue($range, [string]$value) {
    # Forces the value to be stored as text even when it looks like a
    # number/date (e.g. "2024.03.08" or "65"), matching how the rest of
    # the sheet stores such values as plain strings.
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- refresh "想去人数" (F column) counters for existing rows 2-11 ---
    $ws.Range("F2").Value2 = 1851
    $ws.Range("F3").Value2 = 259
    $ws.Range("F4").Value2 = 246
    $ws.Range("F5").Value2 = 8231
    $ws.Range("F6").Value2 = 568
    $ws.Range("F7").Value2 = 589
    $ws.Range("F8").Value2 = 85
    $ws.Range("F10").Value2 = 9199
    $ws.Range("F11").Value2 = 2392

    # --- shift old rows 15,14,13,12 down to 16,15,14,13 (bottom-up so
    #     we never clobber a row before reading its old value) ---

    # row 16 is brand new: copy A15's look (bold/centered/bordered) over
    # before writing its value so the index column stays consistent.
    $ws.Range("A15").Copy()
    $ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

    # old row 15 -> row 16
    $ws.Range("A16").Value2 = 15
    Set-TextValue $ws.Range("B16") $ws.Range("B15").Value2
    Set-TextValue $ws.Range("C16") $ws.Range("C15").Value2
    Set-TextValue $ws.Range("D16") $ws.Range("D15").Value2
    Set-TextValue $ws.Range("E16") $ws.Range("E15").Value2
    $ws.Range("F16").Value2 = 10555
    Set-TextValue $ws.Range("G16") $ws.Range("G15").Value2
    $ws.Range("H16").Value2 = $ws.Range("H15").Value2
    Set-TextValue $ws.Range("I16") $ws.Range("I15").Value2
    Set-TextValue $ws.Range("J16") $ws.Range("J15").Value2

    # old row 14 -> row 15
    $ws.Range("A15").Value2 = 14
    Set-TextValue $ws.Range("B15") $ws.Range("B14").Value2
    Set-TextValue $ws.Range("C15") $ws.Range("C14").Value2
    Set-TextValue $ws.Range("D15") $ws.Range("D14").Value2
    Set-TextValue $ws.Range("E15") $ws.Range("E14").Value2
    $ws.Range("F15").Value2 = 10159
    Set-TextValue $ws.Range("G15") $ws.Range("G14").Value2
    $ws.Range("H15").Value2 = $ws.Range("H14").Value2
    Set-TextValue $ws.Range("I15") $ws.Range("I14").Value2
    Set-TextValue $ws.Range("J15") $ws.Range("J14").Value2

    # old row 13 -> row 14
    $ws.Range("A14").Value2 = 13
    Set-TextValue $ws.Range("B14") $ws.Range("B13").Value2
    Set-TextValue $ws.Range("C14") $ws.Range("C13").Value2
    Set-TextValue $ws.Range("D14") $ws.Range("D13").Value2
    Set-TextValue $ws.Range("E14") $ws.Range("E13").Value2
    $ws.Range("F14").Value2 = 320
    Set-TextValue $ws.Range("G14") $ws.Range("G13").Value2
    $ws.Range("H14").Value2 = $ws.Range("H13").Value2
    Set-TextValue $ws.Range("I14") $ws.Range("I13").Value2
    Set-TextValue $ws.Range("J14") $ws.Range("J13").Value2

    # old row 12 -> row 13
    $ws.Range("A13").Value2 = 12
    Set-TextValue $ws.Range("B13") $ws.Range("B12").Value2
    Set-TextValue $ws.Range("C13") $ws.Range("C12").Value2
    Set-TextValue $ws.Range("D13") $ws.Range("D12").Value2
    Set-TextValue $ws.Range("E13") $ws.Range("E12").Value2
    $ws.Range("F13").Value2 = 28
    Set-TextValue $ws.Range("G13") $ws.Range("G12").Value2
    $ws.Range("H13").Value2 = $ws.Range("H12").Value2
    Set-TextValue $ws.Range("I13") $ws.Range("I12").Value2
    Set-TextValue $ws.Range("J13") $ws.Range("J12").Value2

    # --- brand-new row 12: 2024.03.08 苏州·国风宠物-cosplay展 ---
    $ws.Range("A12").Value2 = 11
    Set-TextValue $ws.Range("B12") "2024.03.08"
    Set-TextValue $ws.Range("C12") "苏州·国风宠物-cosplay展"
    Set-TextValue $ws.Range("D12") "木渎金山南路288号 苏州国际影视娱乐城"
    Set-TextValue $ws.Range("E12") "2024.03.08 09:00-03.10 17:30"
    $ws.Range("F12").Value2 = 5
    Set-TextValue $ws.Range("G12") "65"
    $ws.Range("H12").Value2 = $true
    Set-TextValue $ws.Range("I12") "https://show.bilibili.com/platform/detail.html?id=80635&msource=Msearch_colligation"
    Set-TextValue $ws.Range("J12") "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"
}
